$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the fill color used by the "sub-header" row (old row 10, style with
# theme8 tint 0.8) before it gets removed, so the merged header row can reuse it.
$subHeaderColor = $ws.Range("B10").Interior.Color

# The "Case" sub-header row (old row 10) is removed; Excel shifts rows 11-12
# up to 10-11, and the blank cells to the right of A9 stay on row 9.
$ws.Rows(10).Delete() | Out-Null

# Merge the removed sub-header's column labels into row 9 (the blue header
# row), and add a new "autogen" label under the Action column.
$ws.Range("B9").Value2 = "relationship name"
$ws.Range("C9").Value2 = "doc name"
$ws.Range("D9").Value2 = "desc"
$ws.Range("F9").Value2 = "autogen"

$headerRange = $ws.Range("B9:D9")
$headerRange.Interior.Color = $subHeaderColor
$ws.Range("F9").Interior.Color = $subHeaderColor

# E9 has no value/formatting in the new layout - drop the highlight it
# inherited from the old row 9 (it used to span A9:F9).
$ws.Range("E9").ClearFormats() | Out-Null

# Update the remembered selection to match the new layout.
$ws.Range("E9").Select() | Out-Null
